$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 13
$ws.Range("A3").Value = 42
$ws.Range("A7").Value = 27
$ws.Range("A8").Value = 29
$ws.Range("E8").Value = 473
$ws.Range("A12").Value = 56
$ws.Range("A17").Value = 52
$ws.Range("A19").Value = 48
$ws.Range("E21").Value = 528
$ws.Range("A22").Value = 36
$ws.Range("C22").Value = 'Indiana Government Center Conference Rooms A and B'
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("A23").Value = 51
$ws.Range("C23").Value = 'Indiana Roof Ballroom'
$ws.Range("D23").Value = 4.8
$ws.Range("E23").Value = 284
$ws.Range("A24").Value = 47
$ws.Range("C24").Value = 'Indiana Water Environment Association'
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("A25").Value = 5
$ws.Range("C25").Value = 'Indiana Wesleyan University - Greenwood Education and Conference Center'
$ws.Range("D25").Value = 4.3
$ws.Range("E25").Value = 7
$ws.Range("A26").Value = 7
$ws.Range("C26").Value = 'Indiana Wesleyan University - Indianapolis North Education and Conference Center'
$ws.Range("D26").Value = 4.8
$ws.Range("E26").Value = 12
$ws.Range("A27").Value = 6
$ws.Range("C27").Value = 'Indiana Wesleyan University - Indianapolis West Education and Conference Center'
$ws.Range("D27").Value = 4.3
$ws.Range("E27").Value = 3
$ws.Range("A28").Value = 26
$ws.Range("C28").Value = 'Indiana-Kentucky Conference'
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("A29").Value = 49
$ws.Range("C29").Value = 'Indianapolis EMS Conference Center'
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 1
$ws.Range("A30").Value = 12
$ws.Range("C30").Value = 'Indy West'
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = 2
$ws.Range("A31").Value = 55
$ws.Range("C31").Value = 'IndyFurCon'
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = 3
$ws.Range("A32").Value = 9
$ws.Range("C32").Value = 'Ivy Tech Community College Culinary and Conference Center'
$ws.Range("D32").Value = 4.4
$ws.Range("E32").Value = 84
$ws.Range("A33").Value = 20
$ws.Range("C33").Value = 'JW Marriott Convention Center'
$ws.Range("D33").Value = 4.5
$ws.Range("E33").Value = 13
$ws.Range("A34").Value = 53
$ws.Range("C34").Value = 'Kheprw Institute'
$ws.Range("D34").Value = 4.8
$ws.Range("E34").Value = 32
$ws.Range("A35").Value = 50
$ws.Range("C35").Value = 'MCON'
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("A37").Value = 43
$ws.Range("A39").Value = 28
$ws.Range("A44").Value = 57
$ws.Range("A46").Value = 14
$ws.Range("E46").Value = 845
$ws.Range("A48").Value = 54
$ws.Range("E48").Value = 94
$ws.Range("A52").Value = 24
$ws.Range("A53").Value = 31
$ws.Range("A55").Value = 25
$ws.Range("A56").Value = 11
